$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.433.84"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.13%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.726.17"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +0.07%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'0.9994"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.02%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'243.88"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.30%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'0.9997"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -0.02%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.4911"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +1.80%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.2617"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -1.87%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.06208"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +0.47%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'1.719.81"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -0.36%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.07027"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -2.20%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'15.55"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -0.04%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'4.571"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +1.09%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.6027"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -1.18%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'77.37"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +0.43%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.9995"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -0.04%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'26.440.57"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -0.16%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'0.9995"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Value = "'0.000007198"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +3.64%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'  -1.05%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'1.942.44"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -0.63%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  -0.86%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'8.607"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -1.84%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'5.174"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -1.27%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'137.68"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +0.46%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  -0.43%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'107.11"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +0.42%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  -0.46%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'1.705"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -4.01%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'3.966"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +0.12%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'0.07976"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -0.44%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'3.684"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +0.03%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'0.04529"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +0.60%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'2.600"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -0.59%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'1.000"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +0.76%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'0.6275"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +0.40%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'0.9131"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +0.14%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'1.967"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -4.75%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  +0.56%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  -0.47%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D42").Value = "'99.95"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -3.32%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'5.444"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -3.29%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.3855"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -0.10%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'6.735"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -2.67%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.1158"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -1.77%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.05367"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +0.33%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'30.15"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -0.66%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'7.702"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -0.69%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'  -0.75%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  +0.01%  "
$ws.Range("E51").Style = "Normal"
